$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Spon2"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.54332
$ws.Cells.Item(2, 8).Value = 1.62996
$ws.Cells.Item(2, 9).Value = 0.09277509850694737
$ws.Cells.Item(2, 10).Value = 0.09480543614915297
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 121.928739
$ws.Cells.Item(2, 14).Value = 365.786217
$ws.Cells.Item(2, 15).Value = 0.2282232151508951
$ws.Cells.Item(2, 16).Value = 0.2419720431319445
$ws.Cells.Item(2, 17).Value = 66.24632247347999
$ws.Cells.Item(2, 18).Value = 596.21690226132
$ws.Cells.Item(2, 19).Value = 0.02117343126719654
$ws.Cells.Item(2, 20).Value = 0.02294026508502566

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Spon2"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.54332
$ws.Cells.Item(3, 8).Value = 1.62996
$ws.Cells.Item(3, 9).Value = 0.09277509850694737
$ws.Cells.Item(3, 10).Value = 0.09480543614915297
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 147.91433
$ws.Cells.Item(3, 14).Value = 443.74299
$ws.Cells.Item(3, 15).Value = 0.2768624053389947
$ws.Cells.Item(3, 16).Value = 0.2935413991166814
$ws.Cells.Item(3, 17).Value = 80.3648137756
$ws.Cells.Item(3, 18).Value = 723.2833239803999
$ws.Cells.Item(3, 19).Value = 0.02568593692819562
$ws.Cells.Item(3, 20).Value = 0.02782932037108957

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Spon2"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.54332
$ws.Cells.Item(4, 8).Value = 1.62996
$ws.Cells.Item(4, 9).Value = 0.09277509850694737
$ws.Cells.Item(4, 10).Value = 0.09480543614915297
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 83.50496933333334
$ws.Cells.Item(4, 14).Value = 250.514908
$ws.Cells.Item(4, 15).Value = 0.1563025480180701
$ws.Cells.Item(4, 16).Value = 0.1657186665504434
$ws.Cells.Item(4, 17).Value = 45.36991993818667
$ws.Cells.Item(4, 18).Value = 408.32927944368
$ws.Cells.Item(4, 19).Value = 0.01450098428926332
$ws.Cells.Item(4, 20).Value = 0.01571103046037083

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Spon2"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.54332
$ws.Cells.Item(5, 8).Value = 1.62996
$ws.Cells.Item(5, 9).Value = 0.09277509850694737
$ws.Cells.Item(5, 10).Value = 0.09480543614915297
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 91.06846250000001
$ws.Cells.Item(5, 14).Value = 182.136925
$ws.Cells.Item(5, 15).Value = 0.1704597085236707
$ws.Cells.Item(5, 16).Value = 0.1204857969594293
$ws.Cells.Item(5, 17).Value = 49.47931704550001
$ws.Cells.Item(5, 18).Value = 296.8759022730001
$ws.Cells.Item(5, 19).Value = 0.01581441624974909
$ws.Cells.Item(5, 20).Value = 0.01142270853051699

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Spon2"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.54332
$ws.Cells.Item(6, 8).Value = 1.62996
$ws.Cells.Item(6, 9).Value = 0.09277509850694737
$ws.Cells.Item(6, 10).Value = 0.09480543614915297
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 89.83563
$ws.Cells.Item(6, 14).Value = 269.50689
$ws.Cells.Item(6, 15).Value = 0.1681521229683693
$ws.Cells.Item(6, 16).Value = 0.1782820942415013
$ws.Cells.Item(6, 17).Value = 48.8094944916
$ws.Cells.Item(6, 18).Value = 439.2854504244
$ws.Cells.Item(6, 19).Value = 0.01560032977254279
$ws.Cells.Item(6, 20).Value = 0.01690211170214992

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Spon2"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4.616901666666666
$ws.Cells.Item(7, 8).Value = 13.850705
$ws.Cells.Item(7, 9).Value = 0.7883632241071366
$ws.Cells.Item(7, 10).Value = 0.805616167573593
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 121.928739
$ws.Cells.Item(7, 14).Value = 365.786217
$ws.Cells.Item(7, 15).Value = 0.2282232151508951
$ws.Cells.Item(7, 16).Value = 0.2419720431319445
$ws.Cells.Item(7, 17).Value = 562.9329983036649
$ws.Cells.Item(7, 18).Value = 5066.396984732984
$ws.Cells.Item(7, 19).Value = 0.1799227897124563
$ws.Cells.Item(7, 20).Value = 0.1949365900479093

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Spon2"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.616901666666666
$ws.Cells.Item(8, 8).Value = 13.850705
$ws.Cells.Item(8, 9).Value = 0.7883632241071366
$ws.Cells.Item(8, 10).Value = 0.805616167573593
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 147.91433
$ws.Cells.Item(8, 14).Value = 443.74299
$ws.Cells.Item(8, 15).Value = 0.2768624053389947
$ws.Cells.Item(8, 16).Value = 0.2935413991166814
$ws.Cells.Item(8, 17).Value = 682.9059167008832
$ws.Cells.Item(8, 18).Value = 6146.15325030795
$ws.Cells.Item(8, 19).Value = 0.2182681385071067
$ws.Cells.Item(8, 20).Value = 0.2364816969805714

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Spon2"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.616901666666666
$ws.Cells.Item(9, 8).Value = 13.850705
$ws.Cells.Item(9, 9).Value = 0.7883632241071366
$ws.Cells.Item(9, 10).Value = 0.805616167573593
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 83.50496933333334
$ws.Cells.Item(9, 14).Value = 250.514908
$ws.Cells.Item(9, 15).Value = 0.1563025480180701
$ws.Cells.Item(9, 16).Value = 0.1657186665504434
$ws.Cells.Item(9, 17).Value = 385.5342320900155
$ws.Cells.Item(9, 18).Value = 3469.80808881014
$ws.Cells.Item(9, 19).Value = 0.1232231806916862
$ws.Cells.Item(9, 20).Value = 0.1335056370417744

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Spon2"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.616901666666666
$ws.Cells.Item(10, 8).Value = 13.850705
$ws.Cells.Item(10, 9).Value = 0.7883632241071366
$ws.Cells.Item(10, 10).Value = 0.805616167573593
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 91.06846250000001
$ws.Cells.Item(10, 14).Value = 182.136925
$ws.Cells.Item(10, 15).Value = 0.1704597085236707
$ws.Cells.Item(10, 16).Value = 0.1204857969594293
$ws.Cells.Item(10, 17).Value = 420.4541362970209
$ws.Cells.Item(10, 18).Value = 2522.724817782125
$ws.Cells.Item(10, 19).Value = 0.1343841653920838
$ws.Cells.Item(10, 20).Value = 0.09706530599350552

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Spon2"
$ws.Cells.Item(11, 3).Value = "Itgb1"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.616901666666666
$ws.Cells.Item(11, 8).Value = 13.850705
$ws.Cells.Item(11, 9).Value = 0.7883632241071366
$ws.Cells.Item(11, 10).Value = 0.805616167573593
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 89.83563
$ws.Cells.Item(11, 14).Value = 269.50689
$ws.Cells.Item(11, 15).Value = 0.1681521229683693
$ws.Cells.Item(11, 16).Value = 0.1782820942415013
$ws.Cells.Item(11, 17).Value = 414.7622698730499
$ws.Cells.Item(11, 18).Value = 3732.86042885745
$ws.Cells.Item(11, 19).Value = 0.1325649498038033
$ws.Cells.Item(11, 20).Value = 0.1436269375098324

# Row 12
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "Spon2"
$ws.Cells.Item(12, 3).Value = "Itgb1"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.073119
$ws.Cells.Item(12, 8).Value = 0.219357
$ws.Cells.Item(12, 9).Value = 0.01248550104492653
$ws.Cells.Item(12, 10).Value = 0.01275874012697842
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 121.928739
$ws.Cells.Item(12, 14).Value = 365.786217
$ws.Cells.Item(12, 15).Value = 0.2282232151508951
$ws.Cells.Item(12, 16).Value = 0.2419720431319445
$ws.Cells.Item(12, 17).Value = 8.915307466941
$ws.Cells.Item(12, 18).Value = 80.23776720246899
$ws.Cells.Item(12, 19).Value = 0.002849481191242994
$ws.Cells.Item(12, 20).Value = 0.003087258416314494

# Row 13
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "Spon2"
$ws.Cells.Item(13, 3).Value = "Itgb1"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.073119
$ws.Cells.Item(13, 8).Value = 0.219357
$ws.Cells.Item(13, 9).Value = 0.01248550104492653
$ws.Cells.Item(13, 10).Value = 0.01275874012697842
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 147.91433
$ws.Cells.Item(13, 14).Value = 443.74299
$ws.Cells.Item(13, 15).Value = 0.2768624053389947
$ws.Cells.Item(13, 16).Value = 0.2935413991166814
$ws.Cells.Item(13, 17).Value = 10.81534789527
$ws.Cells.Item(13, 18).Value = 97.33813105742999
$ws.Cells.Item(13, 19).Value = 0.003456765851160892
$ws.Cells.Item(13, 20).Value = 0.003745218427839391

# Row 14
$ws.Cells.Item(14, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 2).Value = "Spon2"
$ws.Cells.Item(14, 3).Value = "Itgb1"
$ws.Cells.Item(14, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.073119
$ws.Cells.Item(14, 8).Value = 0.219357
$ws.Cells.Item(14, 9).Value = 0.01248550104492653
$ws.Cells.Item(14, 10).Value = 0.01275874012697842
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 83.50496933333334
$ws.Cells.Item(14, 14).Value = 250.514908
$ws.Cells.Item(14, 15).Value = 0.1563025480180701
$ws.Cells.Item(14, 16).Value = 0.1657186665504434
$ws.Cells.Item(14, 17).Value = 6.105799852684
$ws.Cells.Item(14, 18).Value = 54.95219867415599
$ws.Cells.Item(14, 19).Value = 0.001951515626604293
$ws.Cells.Item(14, 20).Value = 0.002114361400706498

# Row 15
$ws.Cells.Item(15, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(15, 2).Value = "Spon2"
$ws.Cells.Item(15, 3).Value = "Itgb1"
$ws.Cells.Item(15, 4).Value = "MuSCs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.073119
$ws.Cells.Item(15, 8).Value = 0.219357
$ws.Cells.Item(15, 9).Value = 0.01248550104492653
$ws.Cells.Item(15, 10).Value = 0.01275874012697842
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 91.06846250000001
$ws.Cells.Item(15, 14).Value = 182.136925
$ws.Cells.Item(15, 15).Value = 0.1704597085236707
$ws.Cells.Item(15, 16).Value = 0.1204857969594293
$ws.Cells.Item(15, 17).Value = 6.658834909537501
$ws.Cells.Item(15, 18).Value = 39.953009457225
$ws.Cells.Item(15, 19).Value = 0.002128274868890163
$ws.Cells.Item(15, 20).Value = 0.001537246972397246

# Row 16
$ws.Cells.Item(16, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 2).Value = "Spon2"
$ws.Cells.Item(16, 3).Value = "Itgb1"
$ws.Cells.Item(16, 4).Value = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.073119
$ws.Cells.Item(16, 8).Value = 0.219357
$ws.Cells.Item(16, 9).Value = 0.01248550104492653
$ws.Cells.Item(16, 10).Value = 0.01275874012697842
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 89.83563
$ws.Cells.Item(16, 14).Value = 269.50689
$ws.Cells.Item(16, 15).Value = 0.1681521229683693
$ws.Cells.Item(16, 16).Value = 0.1782820942415013
$ws.Cells.Item(16, 17).Value = 6.56869142997
$ws.Cells.Item(16, 18).Value = 59.11822286973
$ws.Cells.Item(16, 19).Value = 0.00209946350702819
$ws.Cells.Item(16, 20).Value = 0.002274654909720791

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Spon2"
$ws.Cells.Item(17, 3).Value = "Itgb1"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.3762535
$ws.Cells.Item(17, 8).Value = 0.752507
$ws.Cells.Item(17, 9).Value = 0.0642475070420447
$ws.Cells.Item(17, 10).Value = 0.0437690215344491
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 121.928739
$ws.Cells.Item(17, 14).Value = 365.786217
$ws.Cells.Item(17, 15).Value = 0.2282232151508951
$ws.Cells.Item(17, 16).Value = 0.2419720431319445
$ws.Cells.Item(17, 17).Value = 45.8761147993365
$ws.Cells.Item(17, 18).Value = 275.256688796019
$ws.Cells.Item(17, 19).Value = 0.01466277262256521
$ws.Cells.Item(17, 20).Value = 0.01059087956657673

# Row 18
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Spon2"
$ws.Cells.Item(18, 3).Value = "Itgb1"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.3762535
$ws.Cells.Item(18, 8).Value = 0.752507
$ws.Cells.Item(18, 9).Value = 0.0642475070420447
$ws.Cells.Item(18, 10).Value = 0.0437690215344491
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 147.91433
$ws.Cells.Item(18, 14).Value = 443.74299
$ws.Cells.Item(18, 15).Value = 0.2768624053389947
$ws.Cells.Item(18, 16).Value = 0.2935413991166814
$ws.Cells.Item(18, 17).Value = 55.65328436265499
$ws.Cells.Item(18, 18).Value = 333.91970617593
$ws.Cells.Item(18, 19).Value = 0.0177877193366945
$ws.Cells.Item(18, 20).Value = 0.01284801981919035

# Row 19
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Spon2"
$ws.Cells.Item(19, 3).Value = "Itgb1"
$ws.Cells.Item(19, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0.3762535
$ws.Cells.Item(19, 8).Value = 0.752507
$ws.Cells.Item(19, 9).Value = 0.0642475070420447
$ws.Cells.Item(19, 10).Value = 0.0437690215344491
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 83.50496933333334
$ws.Cells.Item(19, 14).Value = 250.514908
$ws.Cells.Item(19, 15).Value = 0.1563025480180701
$ws.Cells.Item(19, 16).Value = 0.1657186665504434
$ws.Cells.Item(19, 17).Value = 31.41903697905934
$ws.Cells.Item(19, 18).Value = 188.514221874356
$ws.Cells.Item(19, 19).Value = 0.01004204905448049
$ws.Cells.Item(19, 20).Value = 0.007253343884906546

# Row 20
$ws.Cells.Item(20, 1).Value = "MuSCs"
$ws.Cells.Item(20, 2).Value = "Spon2"
$ws.Cells.Item(20, 3).Value = "Itgb1"
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 0.3762535
$ws.Cells.Item(20, 8).Value = 0.752507
$ws.Cells.Item(20, 9).Value = 0.0642475070420447
$ws.Cells.Item(20, 10).Value = 0.0437690215344491
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 91.06846250000001
$ws.Cells.Item(20, 14).Value = 182.136925
$ws.Cells.Item(20, 15).Value = 0.1704597085236707
$ws.Cells.Item(20, 16).Value = 0.1204857969594293
$ws.Cells.Item(20, 17).Value = 34.26482775524376
$ws.Cells.Item(20, 18).Value = 137.059311020975
$ws.Cells.Item(20, 19).Value = 0.01095161132375942
$ws.Cells.Item(20, 20).Value = 0.005273545441712524

# Row 21
$ws.Cells.Item(21, 1).Value = "MuSCs"
$ws.Cells.Item(21, 2).Value = "Spon2"
$ws.Cells.Item(21, 3).Value = "Itgb1"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 0.3762535
$ws.Cells.Item(21, 8).Value = 0.752507
$ws.Cells.Item(21, 9).Value = 0.0642475070420447
$ws.Cells.Item(21, 10).Value = 0.0437690215344491
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 89.83563
$ws.Cells.Item(21, 14).Value = 269.50689
$ws.Cells.Item(21, 15).Value = 0.1681521229683693
$ws.Cells.Item(21, 16).Value = 0.1782820942415013
$ws.Cells.Item(21, 17).Value = 33.800970212205
$ws.Cells.Item(21, 18).Value = 202.80582127323
$ws.Cells.Item(21, 19).Value = 0.01080335470454507
$ws.Cells.Item(21, 20).Value = 0.007803232822062954

# Row 22
$ws.Cells.Item(22, 1).Value = "Resolving-Mac"
$ws.Cells.Item(22, 2).Value = "Spon2"
$ws.Cells.Item(22, 3).Value = "Itgb1"
$ws.Cells.Item(22, 4).Value = "ECs"
$ws.Cells.Item(22, 5).Value = 1
$ws.Cells.Item(22, 6).Value = 0.3333333333333333
$ws.Cells.Item(22, 7).Value = 0.2467186666666667
$ws.Cells.Item(22, 8).Value = 0.740156
$ws.Cells.Item(22, 9).Value = 0.04212866929894484
$ws.Cells.Item(22, 10).Value = 0.04305063461582644
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 121.928739
$ws.Cells.Item(22, 14).Value = 365.786217
$ws.Cells.Item(22, 15).Value = 0.2282232151508951
$ws.Cells.Item(22, 16).Value = 0.2419720431319445
$ws.Cells.Item(22, 17).Value = 30.082095914428
$ws.Cells.Item(22, 18).Value = 270.738863229852
$ws.Cells.Item(22, 19).Value = 0.009614740357433997
$ws.Cells.Item(22, 20).Value = 0.01041705001611834

# Row 23
$ws.Cells.Item(23, 1).Value = "Resolving-Mac"
$ws.Cells.Item(23, 2).Value = "Spon2"
$ws.Cells.Item(23, 3).Value = "Itgb1"
$ws.Cells.Item(23, 4).Value = "FAPs"
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = 0.3333333333333333
$ws.Cells.Item(23, 7).Value = 0.2467186666666667
$ws.Cells.Item(23, 8).Value = 0.740156
$ws.Cells.Item(23, 9).Value = 0.04212866929894484
$ws.Cells.Item(23, 10).Value = 0.04305063461582644
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = 147.91433
$ws.Cells.Item(23, 14).Value = 443.74299
$ws.Cells.Item(23, 15).Value = 0.2768624053389947
$ws.Cells.Item(23, 16).Value = 0.2935413991166814
$ws.Cells.Item(23, 17).Value = 36.49322627849333
$ws.Cells.Item(23, 18).Value = 328.43903650644
$ws.Cells.Item(23, 19).Value = 0.01166384471583693
$ws.Cells.Item(23, 20).Value = 0.01263714351799073

# Row 24
$ws.Cells.Item(24, 1).Value = "Resolving-Mac"
$ws.Cells.Item(24, 2).Value = "Spon2"
$ws.Cells.Item(24, 3).Value = "Itgb1"
$ws.Cells.Item(24, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = 0.3333333333333333
$ws.Cells.Item(24, 7).Value = 0.2467186666666667
$ws.Cells.Item(24, 8).Value = 0.740156
$ws.Cells.Item(24, 9).Value = 0.04212866929894484
$ws.Cells.Item(24, 10).Value = 0.04305063461582644
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(24, 13).Value = 83.50496933333334
$ws.Cells.Item(24, 14).Value = 250.514908
$ws.Cells.Item(24, 15).Value = 0.1563025480180701
$ws.Cells.Item(24, 16).Value = 0.1657186665504434
$ws.Cells.Item(24, 17).Value = 20.60223469396089
$ws.Cells.Item(24, 18).Value = 185.420112245648
$ws.Cells.Item(24, 19).Value = 0.00658481835603572
$ws.Cells.Item(24, 20).Value = 0.007134293762685117

# Row 25
$ws.Cells.Item(25, 1).Value = "Resolving-Mac"
$ws.Cells.Item(25, 2).Value = "Spon2"
$ws.Cells.Item(25, 3).Value = "Itgb1"
$ws.Cells.Item(25, 4).Value = "MuSCs"
$ws.Cells.Item(25, 5).Value = 1
$ws.Cells.Item(25, 6).Value = 0.3333333333333333
$ws.Cells.Item(25, 7).Value = 0.2467186666666667
$ws.Cells.Item(25, 8).Value = 0.740156
$ws.Cells.Item(25, 9).Value = 0.04212866929894484
$ws.Cells.Item(25, 10).Value = 0.04305063461582644
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 91.06846250000001
$ws.Cells.Item(25, 14).Value = 182.136925
$ws.Cells.Item(25, 15).Value = 0.1704597085236707
$ws.Cells.Item(25, 16).Value = 0.1204857969594293
$ws.Cells.Item(25, 17).Value = 22.46828964338334
$ws.Cells.Item(25, 18).Value = 134.8097378603
$ws.Cells.Item(25, 19).Value = 0.007181240689188254
$ws.Cells.Item(25, 20).Value = 0.005186990021297045

# Row 26
$ws.Cells.Item(26, 1).Value = "Resolving-Mac"
$ws.Cells.Item(26, 2).Value = "Spon2"
$ws.Cells.Item(26, 3).Value = "Itgb1"
$ws.Cells.Item(26, 4).Value = "Resolving-Mac"
$ws.Cells.Item(26, 5).Value = 1
$ws.Cells.Item(26, 6).Value = 0.3333333333333333
$ws.Cells.Item(26, 7).Value = 0.2467186666666667
$ws.Cells.Item(26, 8).Value = 0.740156
$ws.Cells.Item(26, 9).Value = 0.04212866929894484
$ws.Cells.Item(26, 10).Value = 0.04305063461582644
$ws.Cells.Item(26, 11).Value = 3
$ws.Cells.Item(26, 12).Value = 1
$ws.Cells.Item(26, 13).Value = 89.83563
$ws.Cells.Item(26, 14).Value = 269.50689
$ws.Cells.Item(26, 15).Value = 0.1681521229683693
$ws.Cells.Item(26, 16).Value = 0.1782820942415013
$ws.Cells.Item(26, 17).Value = 22.16412685276
$ws.Cells.Item(26, 18).Value = 199.47714167484
$ws.Cells.Item(26, 19).Value = 0.007084025180449937
$ws.Cells.Item(26, 20).Value = 0.007675157297735207
